$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, pushing existing rows 74-81 down to 75-82
$ws.Rows.Item(74).EntireRow.Insert()

# Populate the new row 74 with the new record
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 44918
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100103
$ws.Range("H74").Value = "Frutos de hueso (carozo)"
$ws.Range("I74").Value = 100103003
$ws.Range("J74").Value = "Damasco"
$ws.Range("K74").Value = "Modesto"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 200
$ws.Range("N74").Value = 20000
$ws.Range("O74").Value = 20000
$ws.Range("P74").Value = 20000
$ws.Range("Q74").Value = "$/bandeja 18 kilos"
$ws.Range("R74").Value = "Región de O'Higgins"
$ws.Range("S74").Value = 1111
$ws.Range("T74").Value = 18

Write-Output "done"
